$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

# Copy formatting from the row above (row 91) so the new row matches
# the existing style (bold/border on column A, date format on column E).
$ws.Range("A91:V91").Copy()
$ws.Range("A92:V92").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 91
$ws.Cells.Item($row, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item($row, 3).Value = "premijer-liga-bih"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45262.66666666666
$ws.Cells.Item($row, 6).Value = "Sloga Doboj"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Borac Banja Luka"
$ws.Cells.Item($row, 9).Value = 4
$ws.Cells.Item($row, 10).Value = 3.05
$ws.Cells.Item($row, 11).Value = "01/12/2023 04:12"
$ws.Cells.Item($row, 12).Value = 8.42
$ws.Cells.Item($row, 13).Value = "02/12/2023 15:59"
$ws.Cells.Item($row, 14).Value = 2.97
$ws.Cells.Item($row, 15).Value = "01/12/2023 04:12"
$ws.Cells.Item($row, 16).Value = 4.4
$ws.Cells.Item($row, 17).Value = "02/12/2023 15:59"
$ws.Cells.Item($row, 18).Value = 2.27
$ws.Cells.Item($row, 19).Value = "01/12/2023 04:12"
$ws.Cells.Item($row, 20).Value = 1.26
$ws.Cells.Item($row, 21).Value = "02/12/2023 15:59"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/sloga-doboj-borac-banja-luka/jVGZWjh4/"
